# Test fixture update for "fuzzy worksheet name matching":
# rename the "Root" sheet to "Main root" and make it the active sheet
# (moving the active tab away from "One to many rows").

$wb = $excel.ActiveWorkbook

$root   = $wb.Worksheets.Item(1)
$nodes  = $wb.Worksheets.Item(2)
$leaves = $wb.Worksheets.Item(4)
$many   = $wb.Worksheets.Item(5)

# 1. Rename the first sheet and make it the active / selected tab.
$root.Name = "Main root"
$root.Activate()

# 2. Re-apply the auto-filter scoped names so each filtered sheet keeps a
#    second "_0" suffixed _FilterDatabase entry alongside the original.
$root.Names.Add("_xlnm._FilterDatabase_0", "='Main root'!`$A`$1:`$B`$2")
$nodes.Names.Add("_xlnm._FilterDatabase_0", "=Nodes!`$A`$1:`$D`$4")
$leaves.Names.Add("_xlnm._FilterDatabase_0", "=Leaves!`$A`$1:`$F`$7")
$many.Names.Add("_xlnm._FilterDatabase_0", "='One to many rows'!`$A`$1:`$A`$13")

# 3. Small column-width tweak on the "Leaves" sheet (column F).
$leaves.Columns.Item(6).ColumnWidth = 23.333333333333332
